$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "25.951.61"
$ws.Cells.Item(2, 5).Value = "  -0.22%  "

$ws.Cells.Item(3, 4).Value = "1.619.08"
$ws.Cells.Item(3, 5).Value = "  -0.93%  "

$ws.Cells.Item(4, 5).Value = "  -0.11%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "212.20"
$ws.Cells.Item(5, 5).Value = "  -0.99%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.488"
$ws.Cells.Item(7, 5).Value = "  -3.20%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.0623"
$ws.Cells.Item(8, 5).Value = "  -0.35%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.247"
$ws.Cells.Item(9, 5).Value = "  -1.50%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "18.27"
$ws.Cells.Item(10, 5).Value = "  -1.76%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0791"
$ws.Cells.Item(11, 5).Value = "  -0.35%  "

$ws.Cells.Item(12, 4).Value = "1.845.18"
$ws.Cells.Item(12, 5).Value = "  -0.88%  "

$ws.Cells.Item(13, 4).Value = "1.623.48"
$ws.Cells.Item(13, 5).Value = "  -3.01%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.13"
$ws.Cells.Item(14, 5).Value = "  -1.72%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.521"
$ws.Cells.Item(15, 5).Value = "  -1.84%  "

$ws.Cells.Item(16, 4).Value = "25.966.24"
$ws.Cells.Item(16, 5).Value = "  -0.21%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "61.69"
$ws.Cells.Item(17, 5).Value = "  -0.57%  "

$ws.Cells.Item(18, 4).Value = "0.0₃0735"

$ws.Cells.Item(19, 5).Value = "  -0.08%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "191.56"
$ws.Cells.Item(20, 5).Value = "  +0.44%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "4.24"
$ws.Cells.Item(21, 5).Value = "  -0.73%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "9.50"
$ws.Cells.Item(22, 5).Value = "  -0.95%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "6.02"
$ws.Cells.Item(23, 5).Value = "  -1.97%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "0.131"
$ws.Cells.Item(24, 5).Value = "  -0.30%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "144.08"
$ws.Cells.Item(25, 5).Value = "  +0.41%  "

$ws.Cells.Item(26, 5).Value = "  -0.11%  "

$ws.Cells.Item(27, 5).Value = "  -3.19%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "6.64"
$ws.Cells.Item(28, 5).Value = "  -1.83%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "15.20"
$ws.Cells.Item(29, 5).Value = "  -0.40%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.23"
$ws.Cells.Item(30, 5).Value = "  -0.98%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.0476"
$ws.Cells.Item(31, 5).Value = "  -1.61%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.11"
$ws.Cells.Item(32, 5).Value = "  -1.52%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "3.08"
$ws.Cells.Item(33, 5).Value = "  -2.78%  "

$ws.Cells.Item(34, 5).Value = "  -0.87%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.40"
$ws.Cells.Item(35, 5).Value = "  -1.36%  "

$ws.Cells.Item(36, 4).Value = "1.128.63"
$ws.Cells.Item(36, 5).Value = "  -0.16%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.834"
$ws.Cells.Item(37, 5).Value = "  -4.83%  "

$ws.Cells.Item(38, 5).Value = "  -1.71%  "

$ws.Cells.Item(39, 5).Value = "  -1.30%  "

$ws.Cells.Item(40, 5).Value = "  -1.22%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "97.72"
$ws.Cells.Item(41, 5).Value = "  -1.24%  "

$ws.Cells.Item(42, 4).Value = "1.756.53"
$ws.Cells.Item(42, 5).Value = "  -0.74%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.755"
$ws.Cells.Item(43, 5).Value = "  -4.17%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "5.12"
$ws.Cells.Item(44, 5).Value = "  -3.51%  "

$ws.Cells.Item(45, 5).Value = "  -0.95%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "1.49"
$ws.Cells.Item(46, 5).Value = "  +0.18%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "53.96"
$ws.Cells.Item(47, 5).Value = "  -2.60%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.0514"
$ws.Cells.Item(48, 5).Value = "  -1.79%  "

$ws.Cells.Item(49, 5).Value = "  -0.71%  "

$ws.Cells.Item(50, 2).Value = "EnergySwap"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "7.47"
$ws.Cells.Item(50, 5).Value = "  -0.50%  "

$ws.Cells.Item(51, 2).Value = "USDD"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.00"
$ws.Cells.Item(51, 5).Value = "  -0.24%  "
